$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text that looks numeric (e.g. '47.00') but must
# stay literal text (matches the source inlineStr cells, incl. trailing zeros
# and thousand-dot formatting). Force text entry, then restore the default
# 'Normal' style so no stray number-format/style id is left on the cell.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '95.341.13'
$ws.Range("E2").Value = '  -1.64%  '
Set-TextValue "D3" '3.616.90'
$ws.Range("E3").Value = '  -2.21%  '
$ws.Range("E4").Value = '  -0.08%  '
Set-TextValue "D5" '2.35'
$ws.Range("E5").Value = '  +21.61%  '
Set-TextValue "D6" '226.87'
$ws.Range("E6").Value = '  -4.40%  '
Set-TextValue "D7" '637.52'
$ws.Range("E7").Value = '  -2.82%  '
Set-TextValue "D8" '0.413'
$ws.Range("E8").Value = '  -3.20%  '
Set-TextValue "D9" '1.10'
$ws.Range("E9").Value = '  +2.88%  '
$ws.Range("E10").Value = '  -0.02%  '
Set-TextValue "D11" '3.617.89'
$ws.Range("E11").Value = '  -2.07%  '
Set-TextValue "D12" '47.00'
$ws.Range("E12").Value = '  +6.38%  '
$ws.Range("E13").Value = '  -0.86%  '
Set-TextValue "D14" '0.0000292'
$ws.Range("E14").Value = '  -1.26%  '
Set-TextValue "D15" '6.49'
$ws.Range("E15").Value = '  -4.05%  '
Set-TextValue "D16" '4.291.52'
$ws.Range("E16").Value = '  -2.23%  '
Set-TextValue "D17" '94.952.55'
$ws.Range("E17").Value = '  -1.75%  '
Set-TextValue "D18" '8.79'
$ws.Range("E18").Value = '  -1.89%  '
Set-TextValue "D19" '3.616.12'
$ws.Range("E19").Value = '  -2.15%  '
Set-TextValue "D20" '19.64'
$ws.Range("E20").Value = '  +5.07%  '
Set-TextValue "D21" '12.82'
$ws.Range("E21").Value = '  -1.31%  '
Set-TextValue "D22" '0.517'
$ws.Range("E22").Value = '  +1.46%  '
Set-TextValue "D23" '512.13'
$ws.Range("E23").Value = '  -1.95%  '
Set-TextValue "D24" '3.25'
$ws.Range("E24").Value = '  -4.72%  '
Set-TextValue "D25" '0.244'
$ws.Range("E25").Value = '  +24.16%  '
Set-TextValue "D26" '120.25'
$ws.Range("E26").Value = '  +18.01%  '
Set-TextValue "D27" '0.0000203'
$ws.Range("E27").Value = '  -3.29%  '
Set-TextValue "D28" '6.74'
$ws.Range("E28").Value = '  -2.53%  '
Set-TextValue "D29" '12.62'
$ws.Range("E29").Value = '  -5.97%  '
Set-TextValue "D30" '12.67'
$ws.Range("E30").Value = '  +3.26%  '
Set-TextValue "D31" '2.92'
$ws.Range("E31").Value = '  -2.77%  '
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("E34").Value = '  -5.44%  '
Set-TextValue "D35" '1.77'
$ws.Range("E35").Value = '  -4.77%  '
Set-TextValue "D36" '31.81'
$ws.Range("E36").Value = '  -1.45%  '
Set-TextValue "D37" '0.587'
$ws.Range("E37").Value = '  -1.78%  '
$ws.Range("E38").Value = '  -0.02%  '
Set-TextValue "D39" '595.12'
$ws.Range("E39").Value = '  -8.05%  '
Set-TextValue "D40" '8.33'
$ws.Range("E40").Value = '  -5.66%  '
Set-TextValue "D41" '6.82'
$ws.Range("E41").Value = '  -0.94%  '
Set-TextValue "D42" '0.486'
$ws.Range("E42").Value = '  +8.46%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D43" '40.05'
$ws.Range("E43").Value = '  -1.04%  '
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D44" '0.158'
$ws.Range("E44").Value = '  -1.07%  '
Set-TextValue "D45" '0.0479'
$ws.Range("E45").Value = '  +4.43%  '
Set-TextValue "D47" '0.922'
$ws.Range("E47").Value = '  -3.67%  '
Set-TextValue "D48" '23.45'
$ws.Range("E48").Value = '  -0.79%  '
Set-TextValue "D49" '8.59'
$ws.Range("E49").Value = '  +0.51%  '
Set-TextValue "D50" '2.21'
$ws.Range("E50").Value = '  -3.35%  '
Set-TextValue "D51" '53.89'
$ws.Range("E51").Value = '  +0.01%  '
